# StaffID.xlsx edit: NATHAN's row (row 4) last name is corrected from the
# placeholder "TEST" to "DANSKIN" (matching row 3's surname). Once "TEST" is
# no longer referenced anywhere in the sheet, the shared-string table drops
# it on save and every later entry (e.g. "2222") shifts down one slot - that
# re-indexing is handled automatically by the engine, we just need to set
# the cell's display value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Staff")

$ws.Range("B4").Value = "DANSKIN"

# Note: the source commit also widens the saved multi-area selection from
# E6,E4,E5 to E6,E3,E4,E5 (same active cell E6). This COM surface only
# supports a single contiguous Range for Select()/Selection, so a
# non-contiguous sqref can't be reproduced faithfully; leaving the
# selection as-is keeps the one attribute (activeCell="E6") that is
# reachable and otherwise unchanged by the edit.
